$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 58, shifting rows 58:84 down to 59:85
$ws.Rows.Item(58).Insert()

# Populate the new row 58 with the new data record
$ws.Cells.Item(58, 1).Value = 9
$ws.Cells.Item(58, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(58, 3).Value = "Metropolitana"
$ws.Cells.Item(58, 4).Value = (Get-Date -Year 2022 -Month 1 -Day 11 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(58, 4).NumberFormat = $ws.Cells.Item(59, 4).NumberFormat
$ws.Cells.Item(58, 5).Value = 13
$ws.Cells.Item(58, 6).Value = "Fruta"
$ws.Cells.Item(58, 7).Value = 100103
$ws.Cells.Item(58, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(58, 9).Value = 100103003
$ws.Cells.Item(58, 10).Value = "Damasco"
$ws.Cells.Item(58, 11).Value = "Patterson"
$ws.Cells.Item(58, 12).Value = "Primera"
$ws.Cells.Item(58, 13).Value = 580
$ws.Cells.Item(58, 14).Value = 14000
$ws.Cells.Item(58, 15).Value = 15000
$ws.Cells.Item(58, 16).Value = 14517
$ws.Cells.Item(58, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(58, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(58, 19).Value = 806
$ws.Cells.Item(58, 20).Value = 18
